# chore: adapt column header formatting to respective input file names (#7)
#
# The sheet compares two AHB "Formatversion" inputs (FV2210 = old, FV2304 =
# new). Column headers used generic "_old"/"_new" suffixes; rename them to
# the concrete format-version they represent, wrap the header block in an
# Excel Table (so the generated suffix is easy to discover/filter on), and
# freeze the header row for easier scrolling through the 64 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: "<Spalte>_old" -> "<Spalte>_FV2210",
#        "<Spalte>_new" -> "<Spalte>_FV2304". "diff" (col K) is unchanged. ---
$headers = @(
    "Segmentname_FV2210", "Segmentgruppe_FV2210", "Segment_FV2210", "Datenelement_FV2210", "Segment ID_FV2210",
    "Code_FV2210", "Qualifier_FV2210", "Beschreibung_FV2210", "Bedingungsausdruck_FV2210", "Bedingung_FV2210",
    "diff",
    "Segmentname_FV2304", "Segmentgruppe_FV2304", "Segment_FV2304", "Datenelement_FV2304", "Segment ID_FV2304",
    "Code_FV2304", "Qualifier_FV2304", "Beschreibung_FV2304", "Bedingungsausdruck_FV2304", "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the header + data block into a proper Excel Table (Table1) so
#        the new header names are exposed as structured-reference column
#        names, matching the renamed headers above. ---
$lastRow = 65
$lastCol = 21
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1, $null)
$tbl.Name = "Table1"

# --- 3. Freeze the header row (row 1) so it stays visible while scrolling. ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()
